$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G3").Value = 2.05
$ws.Range("I3").Value = 3.8
$ws.Range("O3").Value = 1.4
$ws.Range("P3").Value = 2.75
$ws.Range("Y3").Value = 9
$ws.Range("G7").Value = 1.7
$ws.Range("H7").Value = 3.4
$ws.Range("I7").Value = 4.9
$ws.Range("J7").Value = 2.22
$ws.Range("K7").Value = 2.12
$ws.Range("L7").Value = 5
$ws.Range("M7").Value = 1.03
$ws.Range("N7").Value = 6.65
$ws.Range("O7").Value = 1.33
$ws.Range("P7").Value = 2.8
$ws.Range("Q7").Value = 1.98
$ws.Range("S7").Value = 1.39
$ws.Range("T7").Value = 2.57
$ws.Range("W7").Value = 6
$ws.Range("X7").Value = 7.2
$ws.Range("Z7").Value = 12.5
$ws.Range("AA7").Value = 14.5
$ws.Range("AC7").Value = 8.75
$ws.Range("AD7").Value = 6.8
$ws.Range("AG7").Value = 800
$ws.Range("AH7").Value = 12
$ws.Range("AI7").Value = 28
$ws.Range("AJ7").Value = 16
$ws.Range("AL7").Value = 50
$ws.Range("AM7").Value = 55
$ws.Range("AN7").Value = 3.45
$ws.Range("AO7").Value = 8
$ws.Range("AQ7").Value = 27
$ws.Range("AT7").Value = 2.52
$ws.Range("AU7").Value = 7.5
$ws.Range("AW7").Value = 6.4
$ws.Range("AX7").Value = 28
